{"js": "// Update the date line and the division problems in the table to the\n// new values (per commit \"Update master to output generated at 503736d\").\nconst replacements = [\n  [\"2025-02-20 Thursday\", \"2025-02-21 Friday\"],\n  [\"75\u00f74=\", \"26\u00f78=\"],\n  [\"39\u00f77=\", \"52\u00f73=\"],\n  [\"49\u00f73=\", \"12\u00f76=\"],\n  [\"10\u00f77=\", \"20\u00f79=\"],\n  [\"17\u00f78=\", \"51\u00f78=\"],\n  [\"97\u00f74=\", \"12\u00f79=\"],\n  [\"60\u00f72=\", \"70\u00f75=\"],\n  [\"24\u00f74=\", \"24\u00f75=\"],\n  [\"84\u00f75=\", \"67\u00f79=\"],\n  [\"64\u00f73=\", \"98\u00f79=\"],\n  [\"98\u00f75=\", \"34\u00f75=\"],\n  [\"48\u00f73=\", \"56\u00f77=\"],\n  [\"27\u00f75=\", \"58\u00f78=\"],\n  [\"50\u00f77=\", \"65\u00f73=\"],\n  [\"69\u00f73=\", \"23\u00f77=\"],\n  [\"38\u00f76=\", \"39\u00f74=\"],\n  [\"26\u00f77=\", \"33\u00f75=\"],\n  [\"24\u00f72=\", \"93\u00f73=\"],\n  [\"65\u00f76=\", \"19\u00f75=\"],\n  [\"30\u00f77=\", \"59\u00f72=\"],\n  [\"51\u00f75=\", \"29\u00f76=\"],\n  [\"44\u00f75=\", \"32\u00f74=\"],\n  [\"44\u00f78=\", \"42\u00f74=\"],\n  [\"71\u00f76=\", \"84\u00f78=\"],\n  [\"36\u00f77=\", \"93\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the division problems in the table to the\n# new values (per commit \"Update master to output generated at 503736d\").\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-20 Thursday\", \"2025-02-21 Friday\"),\n    @(\"75\u00f74=\", \"26\u00f78=\"),\n    @(\"39\u00f77=\", \"52\u00f73=\"),\n    @(\"49\u00f73=\", \"12\u00f76=\"),\n    @(\"10\u00f77=\", \"20\u00f79=\"),\n    @(\"17\u00f78=\", \"51\u00f78=\"),\n    @(\"97\u00f74=\", \"12\u00f79=\"),\n    @(\"60\u00f72=\", \"70\u00f75=\"),\n    @(\"24\u00f74=\", \"24\u00f75=\"),\n    @(\"84\u00f75=\", \"67\u00f79=\"),\n    @(\"64\u00f73=\", \"98\u00f79=\"),\n    @(\"98\u00f75=\", \"34\u00f75=\"),\n    @(\"48\u00f73=\", \"56\u00f77=\"),\n    @(\"27\u00f75=\", \"58\u00f78=\"),\n    @(\"50\u00f77=\", \"65\u00f73=\"),\n    @(\"69\u00f73=\", \"23\u00f77=\"),\n    @(\"38\u00f76=\", \"39\u00f74=\"),\n    @(\"26\u00f77=\", \"33\u00f75=\"),\n    @(\"24\u00f72=\", \"93\u00f73=\"),\n    @(\"65\u00f76=\", \"19\u00f75=\"),\n    @(\"30\u00f77=\", \"59\u00f72=\"),\n    @(\"51\u00f75=\", \"29\u00f76=\"),\n    @(\"44\u00f75=\", \"32\u00f74=\"),\n    @(\"44\u00f78=\", \"42\u00f74=\"),\n    @(\"71\u00f76=\", \"84\u00f78=\"),\n    @(\"36\u00f77=\", \"93\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
